$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 6720843
$ws.Range("E2").Value = 'Cerro Porteno'
$ws.Range("F2").Value = 'Libertad Asuncion'
$ws.Range("H2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 'H'
$ws.Range("L2").Value = 2.375
$ws.Range("M2").Value = 3.2
$ws.Range("N2").Value = 2.7
$ws.Range("O2").Value = 3.75
$ws.Range("P2").Value = 3.3
$ws.Range("Q2").Value = 1.85
$ws.Range("R2").Value = 0.5
$ws.Range("S2").Value = 1.9
$ws.Range("T2").Value = 1.9
$ws.Range("V2").Value = 1.925
$ws.Range("W2").Value = 1.875
$ws.Range("X2").Value = 2.75
$ws.Range("Y2").Value = -1
$ws.Range("AA2").Value = 0.8999999999999999
$ws.Range("AB2").Value = -1
$ws.Range("AC2").Value = -1
$ws.Range("AD2").Value = 0.875
$ws.Range("B3").Value = 6720873
$ws.Range("E3").Value = 'Sportivo Luqueno'
$ws.Range("F3").Value = 'Sportivo Trinidense'
$ws.Range("H3").Value = 2
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 'D'
$ws.Range("L3").Value = 2.625
$ws.Range("M3").Value = 3.1
$ws.Range("N3").Value = 2.5
$ws.Range("O3").Value = 2.3
$ws.Range("P3").Value = 3.1
$ws.Range("Q3").Value = 2.9
$ws.Range("R3").Value = -0.25
$ws.Range("S3").Value = 2.025
$ws.Range("T3").Value = 1.775
$ws.Range("V3").Value = 1.95
$ws.Range("W3").Value = 1.85
$ws.Range("X3").Value = -1
$ws.Range("Y3").Value = 2.1
$ws.Range("AA3").Value = -0.5
$ws.Range("AB3").Value = 0.3875
$ws.Range("AC3").Value = 0.95
$ws.Range("AD3").Value = -1
$ws.Range("B143").Value = 7493433
$ws.Range("E143").Value = 'Sportivo Luqueno'
$ws.Range("F143").Value = 'Nacional Asuncion'
$ws.Range("G143").Value = 1
$ws.Range("H143").Value = 1
$ws.Range("I143").Value = 1
$ws.Range("J143").Value = 1
$ws.Range("K143").Value = 'D'
$ws.Range("L143").Value = 2.75
$ws.Range("M143").Value = 3.2
$ws.Range("N143").Value = 2.4
$ws.Range("O143").Value = 2.75
$ws.Range("P143").Value = 3.1
$ws.Range("Q143").Value = 2.45
$ws.Range("R143").Value = 0.25
$ws.Range("S143").Value = 1.75
$ws.Range("T143").Value = 2.05
$ws.Range("U143").Value = 2.25
$ws.Range("V143").Value = 2
$ws.Range("W143").Value = 1.8
$ws.Range("X143").Value = -1
$ws.Range("Y143").Value = 2.1
$ws.Range("AA143").Value = 0.375
$ws.Range("AB143").Value = -0.5
$ws.Range("AC143").Value = -0.5
$ws.Range("AD143").Value = 0.4
$ws.Range("B145").Value = 7493312
$ws.Range("E145").Value = 'Cerro Porteno'
$ws.Range("F145").Value = 'Guarani Asuncion'
$ws.Range("G145").Value = 4
$ws.Range("H145").Value = 0
$ws.Range("I145").Value = 3
$ws.Range("J145").Value = 0
$ws.Range("K145").Value = 'H'
$ws.Range("L145").Value = 1.7
$ws.Range("M145").Value = 3.6
$ws.Range("N145").Value = 4.333
$ws.Range("O145").Value = 1.727
$ws.Range("P145").Value = 3.75
$ws.Range("Q145").Value = 4.2
$ws.Range("R145").Value = -0.5
$ws.Range("S145").Value = 1.8
$ws.Range("T145").Value = 2
$ws.Range("U145").Value = 2.75
$ws.Range("V145").Value = 1.875
$ws.Range("W145").Value = 1.925
$ws.Range("X145").Value = 0.7270000000000001
$ws.Range("Y145").Value = -1
$ws.Range("AA145").Value = 0.8
$ws.Range("AB145").Value = -1
$ws.Range("AC145").Value = 0.875
$ws.Range("AD145").Value = -1
$ws.Range("B236").Value = 7609209
$ws.Range("E236").Value = 'Nacional Asuncion'
$ws.Range("F236").Value = 'Sportivo Luqueno'
$ws.Range("G236").Value = 0
$ws.Range("H236").Value = 1
$ws.Range("J236").Value = 0
$ws.Range("L236").Value = 2.9
$ws.Range("M236").Value = 3.1
$ws.Range("N236").Value = 2.5
$ws.Range("P236").Value = 3.25
$ws.Range("Q236").Value = 2.55
$ws.Range("S236").Value = 1.925
$ws.Range("T236").Value = 1.875
$ws.Range("U236").Value = 2.5
$ws.Range("V236").Value = 2
$ws.Range("W236").Value = 1.8
$ws.Range("Z236").Value = 1.55
$ws.Range("AB236").Value = 0.875
$ws.Range("AC236").Value = -1
$ws.Range("AD236").Value = 0.8
$ws.Range("B237").Value = 7609208
$ws.Range("E237").Value = 'Libertad Asuncion'
$ws.Range("F237").Value = 'Cerro Porteno'
$ws.Range("G237").Value = 1
$ws.Range("H237").Value = 3
$ws.Range("J237").Value = 1
$ws.Range("L237").Value = 2.75
$ws.Range("M237").Value = 3
$ws.Range("N237").Value = 2.625
$ws.Range("P237").Value = 2.875
$ws.Range("Q237").Value = 2.8
$ws.Range("S237").Value = 1.875
$ws.Range("T237").Value = 1.925
$ws.Range("U237").Value = 2
$ws.Range("V237").Value = 1.925
$ws.Range("W237").Value = 1.875
$ws.Range("Z237").Value = 1.8
$ws.Range("AB237").Value = 0.925
$ws.Range("AC237").Value = 0.925
$ws.Range("AD237").Value = -1
